$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(432,1).Value = 27354
$ws.Cells.Item(432,2).Value = 45465
$ws.Cells.Item(432,3).Value = 4
$ws.Cells.Item(432,4).Value = "Miami Heretics"
$ws.Cells.Item(432,5).Value = "Rio"
$ws.Cells.Item(432,6).Value = "Vista"
$ws.Cells.Item(432,7).Value = "Rio"
$ws.Cells.Item(432,8).Value = "6 Star"
$ws.Cells.Item(432,9).Value = "Invasion"

$ws.Cells.Item(433,1).Value = 27354
$ws.Cells.Item(433,2).Value = 45465
$ws.Cells.Item(433,3).Value = 4
$ws.Cells.Item(433,4).Value = "Atlanta FaZe"
$ws.Cells.Item(433,5).Value = "Sub Base"
$ws.Cells.Item(433,6).Value = "6 Star"
$ws.Cells.Item(433,7).Value = "Karachi"
$ws.Cells.Item(433,8).Value = "Highrise"
$ws.Cells.Item(433,10).Value = "Highrise"

$ws.Cells.Item(434,1).Value = 27355
$ws.Cells.Item(434,2).Value = 45465
$ws.Cells.Item(434,3).Value = 4
$ws.Cells.Item(434,4).Value = "Seattle Surge"
$ws.Cells.Item(434,5).Value = "Sub Base"
$ws.Cells.Item(434,6).Value = "Vista"
$ws.Cells.Item(434,7).Value = "Rio"
$ws.Cells.Item(434,8).Value = "Karachi"
$ws.Cells.Item(434,9).Value = "Highrise"

$ws.Cells.Item(435,1).Value = 27355
$ws.Cells.Item(435,2).Value = 45465
$ws.Cells.Item(435,3).Value = 4
$ws.Cells.Item(435,4).Value = "Los Angeles Guerrillas"
$ws.Cells.Item(435,5).Value = "6 Star"
$ws.Cells.Item(435,6).Value = "Rio"
$ws.Cells.Item(435,7).Value = "Invasion"
$ws.Cells.Item(435,8).Value = "6 Star"
$ws.Cells.Item(435,10).Value = "Karachi"

$ws.Cells.Item(436,1).Value = 27356
$ws.Cells.Item(436,2).Value = 45465
$ws.Cells.Item(436,3).Value = 4
$ws.Cells.Item(436,4).Value = "Boston Breach"
$ws.Cells.Item(436,5).Value = "Sub Base"
$ws.Cells.Item(436,6).Value = "6 Star"
$ws.Cells.Item(436,7).Value = "Invasion"
$ws.Cells.Item(436,8).Value = "Rio"
$ws.Cells.Item(436,9).Value = "Invasion"

$ws.Cells.Item(437,1).Value = 27356
$ws.Cells.Item(437,2).Value = 45465
$ws.Cells.Item(437,3).Value = 4
$ws.Cells.Item(437,4).Value = "OpTic Texas"
$ws.Cells.Item(437,5).Value = "Vista"
$ws.Cells.Item(437,6).Value = "Rio"
$ws.Cells.Item(437,7).Value = "Highrise"
$ws.Cells.Item(437,8).Value = "6 Star"
$ws.Cells.Item(437,10).Value = "Karachi"

$ws.Cells.Item(438,1).Value = 27357
$ws.Cells.Item(438,2).Value = 45465
$ws.Cells.Item(438,3).Value = 4
$ws.Cells.Item(438,4).Value = "Minnesota ROKKR"
$ws.Cells.Item(438,5).Value = "Sub Base"
$ws.Cells.Item(438,6).Value = "6 Star"
$ws.Cells.Item(438,7).Value = "Invasion"
$ws.Cells.Item(438,8).Value = "6 Star"
$ws.Cells.Item(438,9).Value = "Highrise"

$ws.Cells.Item(439,1).Value = 27357
$ws.Cells.Item(439,2).Value = 45465
$ws.Cells.Item(439,3).Value = 4
$ws.Cells.Item(439,4).Value = "Las Vegas Legion"
$ws.Cells.Item(439,5).Value = "Vista"
$ws.Cells.Item(439,6).Value = "Rio"
$ws.Cells.Item(439,7).Value = "Highrise"
$ws.Cells.Item(439,8).Value = "Rio"
$ws.Cells.Item(439,10).Value = "Invasion"

$ws.Cells.Item(440,1).Value = 27358
$ws.Cells.Item(440,2).Value = 45466
$ws.Cells.Item(440,3).Value = 4
$ws.Cells.Item(440,4).Value = "OpTic Texas"
$ws.Cells.Item(440,5).Value = "Vista"
$ws.Cells.Item(440,6).Value = "Sub Base"
$ws.Cells.Item(440,7).Value = "Rio"
$ws.Cells.Item(440,8).Value = "Invasion"
$ws.Cells.Item(440,9).Value = "Karachi"

$ws.Cells.Item(441,1).Value = 27358
$ws.Cells.Item(441,2).Value = 45466
$ws.Cells.Item(441,3).Value = 4
$ws.Cells.Item(441,4).Value = "Toronto Ultra"
$ws.Cells.Item(441,5).Value = "Rio"
$ws.Cells.Item(441,6).Value = "Karachi"
$ws.Cells.Item(441,7).Value = "Highrise"
$ws.Cells.Item(441,8).Value = "6 Star"
$ws.Cells.Item(441,10).Value = "Highrise"

$ws.Cells.Item(442,1).Value = 27359
$ws.Cells.Item(442,2).Value = 45466
$ws.Cells.Item(442,3).Value = 4
$ws.Cells.Item(442,4).Value = "Carolina Royal Ravens"
$ws.Cells.Item(442,5).Value = "Rio"
$ws.Cells.Item(442,6).Value = "Karachi"
$ws.Cells.Item(442,7).Value = "6 Star"
$ws.Cells.Item(442,8).Value = "Karachi"
$ws.Cells.Item(442,9).Value = "Highrise"

$ws.Cells.Item(443,1).Value = 27359
$ws.Cells.Item(443,2).Value = 45466
$ws.Cells.Item(443,3).Value = 4
$ws.Cells.Item(443,4).Value = "Boston Breach"
$ws.Cells.Item(443,5).Value = "Sub Base"
$ws.Cells.Item(443,6).Value = "6 Star"
$ws.Cells.Item(443,7).Value = "Highrise"
$ws.Cells.Item(443,8).Value = "Rio"
$ws.Cells.Item(443,10).Value = "Karachi"

$ws.Cells.Item(444,1).Value = 27360
$ws.Cells.Item(444,2).Value = 45466
$ws.Cells.Item(444,3).Value = 4
$ws.Cells.Item(444,4).Value = "Los Angeles Thieves"
$ws.Cells.Item(444,5).Value = "Karachi"
$ws.Cells.Item(444,6).Value = "Vista"
$ws.Cells.Item(444,7).Value = "6 Star"
$ws.Cells.Item(444,8).Value = "Highrise"
$ws.Cells.Item(444,9).Value = "Highrise"

$ws.Cells.Item(445,1).Value = 27360
$ws.Cells.Item(445,2).Value = 45466
$ws.Cells.Item(445,3).Value = 4
$ws.Cells.Item(445,4).Value = "New York Subliners"
$ws.Cells.Item(445,5).Value = "6 Star"
$ws.Cells.Item(445,6).Value = "Sub Base"
$ws.Cells.Item(445,7).Value = "Rio"
$ws.Cells.Item(445,8).Value = "Invasion"
$ws.Cells.Item(445,10).Value = "Invasion"

$ws.Cells.Item(446,1).Value = 91336
$ws.Cells.Item(446,2).Value = 45470
$ws.Cells.Item(446,3).Value = 4
$ws.Cells.Item(446,4).Value = "Seattle Surge"
$ws.Cells.Item(446,5).Value = "Sub Base"
$ws.Cells.Item(446,6).Value = "Karachi"
$ws.Cells.Item(446,7).Value = "Rio"
$ws.Cells.Item(446,8).Value = "6 Star"
$ws.Cells.Item(446,9).Value = "Highrise"

$ws.Cells.Item(447,1).Value = 91336
$ws.Cells.Item(447,2).Value = 45470
$ws.Cells.Item(447,3).Value = 4
$ws.Cells.Item(447,4).Value = "Las Vegas Legion"
$ws.Cells.Item(447,5).Value = "Vista"
$ws.Cells.Item(447,6).Value = "Rio"
$ws.Cells.Item(447,7).Value = "Invasion"
$ws.Cells.Item(447,8).Value = "Karachi"
$ws.Cells.Item(447,10).Value = "Invasion"

$ws.Cells.Item(448,1).Value = 91337
$ws.Cells.Item(448,2).Value = 45470
$ws.Cells.Item(448,3).Value = 4
$ws.Cells.Item(448,4).Value = "Toronto Ultra"
$ws.Cells.Item(448,5).Value = "Rio"
$ws.Cells.Item(448,6).Value = "Karachi"
$ws.Cells.Item(448,7).Value = "Highrise"
$ws.Cells.Item(448,8).Value = "Rio"
$ws.Cells.Item(448,9).Value = "Invasion"

$ws.Cells.Item(449,1).Value = 91337
$ws.Cells.Item(449,2).Value = 45470
$ws.Cells.Item(449,3).Value = 4
$ws.Cells.Item(449,4).Value = "Los Angeles Guerrillas"
$ws.Cells.Item(449,5).Value = "Vista"
$ws.Cells.Item(449,6).Value = "6 Star"
$ws.Cells.Item(449,7).Value = "Karachi"
$ws.Cells.Item(449,8).Value = "6 Star"
$ws.Cells.Item(449,10).Value = "Highrise"

$ws.Cells.Item(450,1).Value = 91338
$ws.Cells.Item(450,2).Value = 45470
$ws.Cells.Item(450,3).Value = 4
$ws.Cells.Item(450,4).Value = "New York Subliners"
$ws.Cells.Item(450,5).Value = "Rio"
$ws.Cells.Item(450,6).Value = "6 Star"
$ws.Cells.Item(450,7).Value = "Rio"
$ws.Cells.Item(450,8).Value = "Highrise"
$ws.Cells.Item(450,9).Value = "Karachi"

$ws.Cells.Item(451,1).Value = 91338
$ws.Cells.Item(451,2).Value = 45470
$ws.Cells.Item(451,3).Value = 4
$ws.Cells.Item(451,4).Value = "Atlanta FaZe"
$ws.Cells.Item(451,5).Value = "Vista"
$ws.Cells.Item(451,6).Value = "Karachi"
$ws.Cells.Item(451,7).Value = "Invasion"
$ws.Cells.Item(451,8).Value = "6 Star"
$ws.Cells.Item(451,10).Value = "Invasion"

$ws.Cells.Item(452,1).Value = 91339
$ws.Cells.Item(452,2).Value = 45470
$ws.Cells.Item(452,3).Value = 4
$ws.Cells.Item(452,4).Value = "Miami Heretics"
$ws.Cells.Item(452,5).Value = "6 Star"
$ws.Cells.Item(452,6).Value = "Rio"
$ws.Cells.Item(452,7).Value = "Rio"
$ws.Cells.Item(452,8).Value = "Karachi"
$ws.Cells.Item(452,9).Value = "Invasion"

$ws.Cells.Item(453,1).Value = 91339
$ws.Cells.Item(453,2).Value = 45470
$ws.Cells.Item(453,3).Value = 4
$ws.Cells.Item(453,4).Value = "Los Angeles Thieves"
$ws.Cells.Item(453,5).Value = "Sub Base"
$ws.Cells.Item(453,6).Value = "Karachi"
$ws.Cells.Item(453,7).Value = "6 Star"
$ws.Cells.Item(453,8).Value = "Invasion"
$ws.Cells.Item(453,10).Value = "Highrise"

$ws.Cells.Item(454,1).Value = 91340
$ws.Cells.Item(454,2).Value = 45471
$ws.Cells.Item(454,3).Value = 4
$ws.Cells.Item(454,4).Value = "Carolina Royal Ravens"
$ws.Cells.Item(454,5).Value = "Vista"
$ws.Cells.Item(454,6).Value = "Karachi"
$ws.Cells.Item(454,7).Value = "6 Star"
$ws.Cells.Item(454,8).Value = "Rio"
$ws.Cells.Item(454,9).Value = "Highrise"

$ws.Cells.Item(455,1).Value = 91340
$ws.Cells.Item(455,2).Value = 45471
$ws.Cells.Item(455,3).Value = 4
$ws.Cells.Item(455,4).Value = "Toronto Ultra"
$ws.Cells.Item(455,5).Value = "Rio"
$ws.Cells.Item(455,6).Value = "6 Star"
$ws.Cells.Item(455,7).Value = "Highrise"
$ws.Cells.Item(455,8).Value = "Karachi"
$ws.Cells.Item(455,10).Value = "Karachi"

$ws.Cells.Item(456,1).Formula = "=A454+1"
$ws.Cells.Item(456,2).Value = 45471
$ws.Cells.Item(456,3).Value = 4
$ws.Cells.Item(456,4).Value = "Las Vegas Legion"
$ws.Cells.Item(456,5).Value = "Rio"
$ws.Cells.Item(456,6).Value = "6 Star"
$ws.Cells.Item(456,7).Value = "Invasion"
$ws.Cells.Item(456,8).Value = "6 Star"
$ws.Cells.Item(456,9).Value = "Karachi"

$ws.Cells.Item(457,1).Formula = "=A455+1"
$ws.Cells.Item(457,2).Value = 45471
$ws.Cells.Item(457,3).Value = 4
$ws.Cells.Item(457,4).Value = "Boston Breach"
$ws.Cells.Item(457,5).Value = "Karachi"
$ws.Cells.Item(457,6).Value = "Vista"
$ws.Cells.Item(457,7).Value = "Highrise"
$ws.Cells.Item(457,8).Value = "Rio"
$ws.Cells.Item(457,10).Value = "Highrise"

$ws.Cells.Item(458,1).Formula = "=A456+1"
$ws.Cells.Item(458,2).Value = 45471
$ws.Cells.Item(458,3).Value = 4
$ws.Cells.Item(458,4).Value = "Minnesota ROKKR"
$ws.Cells.Item(458,5).Value = "6 Star"
$ws.Cells.Item(458,6).Value = "Vista"
$ws.Cells.Item(458,7).Value = "Invasion"
$ws.Cells.Item(458,8).Value = "6 Star"
$ws.Cells.Item(458,9).Value = "Highrise"

$ws.Cells.Item(459,1).Formula = "=A457+1"
$ws.Cells.Item(459,2).Value = 45471
$ws.Cells.Item(459,3).Value = 4
$ws.Cells.Item(459,4).Value = "New York Subliners"
$ws.Cells.Item(459,5).Value = "Rio"
$ws.Cells.Item(459,6).Value = "Karachi"
$ws.Cells.Item(459,7).Value = "Rio"
$ws.Cells.Item(459,8).Value = "Highrise"
$ws.Cells.Item(459,10).Value = "Invasion"

$ws.Cells.Item(460,1).Formula = "=A458+1"
$ws.Cells.Item(460,2).Value = 45471
$ws.Cells.Item(460,3).Value = 4
$ws.Cells.Item(460,4).Value = "Miami Heretics"
$ws.Cells.Item(460,5).Value = "Karachi"
$ws.Cells.Item(460,6).Value = "Sub Base"
$ws.Cells.Item(460,7).Value = "Rio"
$ws.Cells.Item(460,8).Value = "Karachi"
$ws.Cells.Item(460,9).Value = "Invasion"

$ws.Cells.Item(461,1).Formula = "=A459+1"
$ws.Cells.Item(461,2).Value = 45471
$ws.Cells.Item(461,3).Value = 4
$ws.Cells.Item(461,4).Value = "OpTic Texas"
$ws.Cells.Item(461,5).Value = "Vista"
$ws.Cells.Item(461,6).Value = "Rio"
$ws.Cells.Item(461,7).Value = "6 Star"
$ws.Cells.Item(461,8).Value = "Invasion"
$ws.Cells.Item(461,10).Value = "Highrise"

$ws.Cells.Item(462,1).Formula = "=A460+1"
$ws.Cells.Item(462,2).Value = 45471
$ws.Cells.Item(462,3).Value = 4
$ws.Cells.Item(462,4).Value = "Las Vegas Legion"
$ws.Cells.Item(462,5).Value = "Vista"
$ws.Cells.Item(462,6).Value = "Karachi"
$ws.Cells.Item(462,7).Value = "Invasion"
$ws.Cells.Item(462,8).Value = "6 Star"
$ws.Cells.Item(462,9).Value = "Karachi"

$ws.Cells.Item(463,1).Formula = "=A461+1"
$ws.Cells.Item(463,2).Value = 45471
$ws.Cells.Item(463,3).Value = 4
$ws.Cells.Item(463,4).Value = "New York Subliners"
$ws.Cells.Item(463,5).Value = "6 Star"
$ws.Cells.Item(463,6).Value = "Rio"
$ws.Cells.Item(463,7).Value = "Rio"
$ws.Cells.Item(463,8).Value = "Highrise"
$ws.Cells.Item(463,10).Value = "Highrise"

$ws.Cells.Item(464,1).Formula = "=A462+1"
$ws.Cells.Item(464,2).Value = 45472
$ws.Cells.Item(464,3).Value = 4
$ws.Cells.Item(464,4).Value = "Miami Heretics"
$ws.Cells.Item(464,5).Value = "Karachi"
$ws.Cells.Item(464,6).Value = "Vista"
$ws.Cells.Item(464,7).Value = "Rio"
$ws.Cells.Item(464,8).Value = "Karachi"
$ws.Cells.Item(464,9).Value = "Karachi"

$ws.Cells.Item(465,1).Formula = "=A463+1"
$ws.Cells.Item(465,2).Value = 45472
$ws.Cells.Item(465,3).Value = 4
$ws.Cells.Item(465,4).Value = "Toronto Ultra"
$ws.Cells.Item(465,5).Value = "Rio"
$ws.Cells.Item(465,6).Value = "6 Star"
$ws.Cells.Item(465,7).Value = "Highrise"
$ws.Cells.Item(465,8).Value = "6 Star"
$ws.Cells.Item(465,10).Value = "Highrise"

$ws.Cells.Item(466,1).Formula = "=A464+1"
$ws.Cells.Item(466,2).Value = 45472
$ws.Cells.Item(466,3).Value = 4
$ws.Cells.Item(466,4).Value = "Seattle Surge"
$ws.Cells.Item(466,5).Value = "Rio"
$ws.Cells.Item(466,6).Value = "Sub Base"
$ws.Cells.Item(466,7).Value = "Rio"
$ws.Cells.Item(466,8).Value = "Highrise"
$ws.Cells.Item(466,9).Value = "Invasion"

$ws.Cells.Item(467,1).Formula = "=A465+1"
$ws.Cells.Item(467,2).Value = 45472
$ws.Cells.Item(467,3).Value = 4
$ws.Cells.Item(467,4).Value = "Atlanta FaZe"
$ws.Cells.Item(467,5).Value = "Vista"
$ws.Cells.Item(467,6).Value = "Karachi"
$ws.Cells.Item(467,7).Value = "Invasion"
$ws.Cells.Item(467,8).Value = "Karachi"
$ws.Cells.Item(467,10).Value = "Highrise"

$ws.Cells.Item(468,1).Formula = "=A466+1"
$ws.Cells.Item(468,2).Value = 45472
$ws.Cells.Item(468,3).Value = 4
$ws.Cells.Item(468,4).Value = "Los Angeles Thieves"
$ws.Cells.Item(468,5).Value = "Karachi"
$ws.Cells.Item(468,6).Value = "Rio"
$ws.Cells.Item(468,7).Value = "6 Star"
$ws.Cells.Item(468,8).Value = "Rio"
$ws.Cells.Item(468,9).Value = "Karachi"

$ws.Cells.Item(469,1).Formula = "=A467+1"
$ws.Cells.Item(469,2).Value = 45472
$ws.Cells.Item(469,3).Value = 4
$ws.Cells.Item(469,4).Value = "Los Angeles Guerrillas"
$ws.Cells.Item(469,5).Value = "Sub Base"
$ws.Cells.Item(469,6).Value = "6 Star"
$ws.Cells.Item(469,7).Value = "Invasion"
$ws.Cells.Item(469,8).Value = "Highrise"
$ws.Cells.Item(469,10).Value = "Invasion"

$ws.Cells.Item(470,1).Formula = "=A468+1"
$ws.Cells.Item(470,2).Value = 45472
$ws.Cells.Item(470,3).Value = 4
$ws.Cells.Item(470,4).Value = "Seattle Surge"
$ws.Cells.Item(470,5).Value = "Sub Base"
$ws.Cells.Item(470,6).Value = "Vista"
$ws.Cells.Item(470,7).Value = "Rio"
$ws.Cells.Item(470,8).Value = "Invasion"
$ws.Cells.Item(470,9).Value = "Highrise"

$ws.Cells.Item(471,1).Formula = "=A469+1"
$ws.Cells.Item(471,2).Value = 45472
$ws.Cells.Item(471,3).Value = 4
$ws.Cells.Item(471,4).Value = "Toronto Ultra"
$ws.Cells.Item(471,5).Value = "Rio"
$ws.Cells.Item(471,6).Value = "Karachi"
$ws.Cells.Item(471,7).Value = "Highrise"
$ws.Cells.Item(471,8).Value = "6 Star"
$ws.Cells.Item(471,10).Value = "Karachi"

$ws.Cells.Item(472,1).Formula = "=A470+1"
$ws.Cells.Item(472,2).Value = 45472
$ws.Cells.Item(472,3).Value = 4
$ws.Cells.Item(472,4).Value = "Los Angeles Thieves"
$ws.Cells.Item(472,5).Value = "Karachi"
$ws.Cells.Item(472,6).Value = "Rio"
$ws.Cells.Item(472,7).Value = "6 Star"
$ws.Cells.Item(472,8).Value = "Highrise"
$ws.Cells.Item(472,9).Value = "Karachi"

$ws.Cells.Item(473,1).Formula = "=A471+1"
$ws.Cells.Item(473,2).Value = 45472
$ws.Cells.Item(473,3).Value = 4
$ws.Cells.Item(473,4).Value = "New York Subliners"
$ws.Cells.Item(473,5).Value = "6 Star"
$ws.Cells.Item(473,6).Value = "Vista"
$ws.Cells.Item(473,7).Value = "Rio"
$ws.Cells.Item(473,8).Value = "Invasion"
$ws.Cells.Item(473,10).Value = "Invasion"
